$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.740.02"
$ws.Range("E2").Value = "  -0.80%  "

$ws.Range("D3").Value = "2.447.79"
$ws.Range("E3").Value = "  +0.39%  "

$ws.Range("E4").Value = "  -0.06%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "577.25"
$ws.Range("E5").Value = "  -0.49%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "140.98"
$ws.Range("E6").Value = "  -1.30%  "

$ws.Range("E7").Value = "  +0.10%  "

$ws.Range("E8").Value = "  +0.70%  "

$ws.Range("D9").Value = "2.438.79"
$ws.Range("E9").Value = "  +0.14%  "

$ws.Range("E10").Value = "  +2.82%  "

$ws.Range("E11").Value = "  +1.78%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.16"
$ws.Range("E12").Value = "  -0.55%  "

$ws.Range("E13").Value = "  -1.01%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "26.01"
$ws.Range("E14").Value = "  -1.29%  "

$ws.Range("D15").Value = "2.891.55"
$ws.Range("E15").Value = "  +0.80%  "

$ws.Range("E16").Value = "  -0.36%  "

$ws.Range("D17").Value = "61.694.76"
$ws.Range("E17").Value = "  -0.71%  "

$ws.Range("D18").Value = "2.455.21"
$ws.Range("E18").Value = "  +0.74%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.60"
$ws.Range("E19").Value = "  -3.15%  "

$ws.Range("E20").Value = "  +2.40%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "325.05"
$ws.Range("E21").Value = "  -1.66%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.06"
$ws.Range("E22").Value = "  -0.96%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.05"
$ws.Range("E23").Value = "  +1.24%  "

$ws.Range("B24").Value = "Dai"
$ws.Range("C24").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.00"
$ws.Range("E24").Value = "  +0.04%  "

$ws.Range("B25").Value = "SuiNetwork"
$ws.Range("C25").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.93"
$ws.Range("E25").Value = "  -0.87%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "64.88"
$ws.Range("E26").Value = "  -1.10%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.14"
$ws.Range("E27").Value = "  -1.14%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "582.53"
$ws.Range("E28").Value = "  -7.32%  "

$ws.Range("E30").Value = "  +0.23%  "

$ws.Range("D31").Value = "0.0₃0923"
$ws.Range("E31").Value = "  -2.73%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.92"
$ws.Range("E32").Value = "  -0.72%  "

$ws.Range("E33").Value = "  -4.61%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.86"
$ws.Range("E34").Value = "  -0.85%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.133"
$ws.Range("E35").Value = "  -6.00%  "

$ws.Range("E36").Value = "  +0.12%  "

$ws.Range("E37").Value = "  -5.01%  "

$ws.Range("E38").Value = "  -0.49%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "151.79"
$ws.Range("E39").Value = "  +1.88%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.39"
$ws.Range("E40").Value = "  -2.93%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "18.31"
$ws.Range("E41").Value = "  -0.02%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.14"
$ws.Range("E42").Value = "  -1.56%  "

$ws.Range("E43").Value = "  +0.02%  "

$ws.Range("E44").Value = "  -4.09%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "41.71"
$ws.Range("E45").Value = "  -2.45%  "

$ws.Range("E46").Value = "  -3.80%  "

$ws.Range("D47").Value = "0.0₆0295"
$ws.Range("E47").Value = "  +25.80%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "143.09"
$ws.Range("E48").Value = "  +0.24%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.56"
$ws.Range("E49").Value = "  -2.02%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.598"
$ws.Range("E50").Value = "  +0.37%  "

$ws.Range("B51").Value = "Hedera"
$ws.Range("C51").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0510"
$ws.Range("E51").Value = "  -2.19%  "
